$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.640.76"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "2.272.10"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'118.55"
$ws.Range("E5").Value = "  +5.56%  "
$ws.Range("D6").Value = "'268.15"
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("D7").Value = "'0.637"
$ws.Range("E7").Value = "  +2.74%  "
$ws.Range("D8").Value = "'1.01"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.622"
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("D10").Value = "'47.32"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "'0.0945"
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("D12").Value = "'9.43"
$ws.Range("E12").Value = "  +7.13%  "
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "'15.87"
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").Value = "'0.921"
$ws.Range("E15").Value = "  +7.80%  "
$ws.Range("D16").Value = "2.614.23"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "2.270.28"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "43.565.67"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").Value = "'6.93"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").Value = "'72.36"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("E22").Value = "  -4.18%  "
$ws.Range("D23").Value = "'234.74"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("E24").Value = "  +3.25%  "
$ws.Range("D25").Value = "'9.67"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "'12.33"
$ws.Range("E26").Value = "  +9.01%  "
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("D28").Value = "'41.79"
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").Value = "'174.60"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").Value = "'5.74"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("E35").Value = "  +3.13%  "
$ws.Range("D36").Value = "'4.29"
$ws.Range("E36").Value = "  +12.00%  "
$ws.Range("D37").Value = "'0.0381"
$ws.Range("E37").Value = "  +8.44%  "
$ws.Range("D38").Value = "'4.64"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  +3.55%  "
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("D41").Value = "'13.87"
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("D43").Value = "'72.19"
$ws.Range("E43").Value = "  -5.63%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").Value = "'5.74"
$ws.Range("E46").Value = "  -5.76%  "
$ws.Range("D47").Value = "'0.685"
$ws.Range("E47").Value = "  +22.19%  "
$ws.Range("D48").Value = "'74.49"
$ws.Range("E48").Value = "  +36.86%  "
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").Value = "'103.37"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("E51").Value = "  -0.22%  "
